# Refined metadata to be additional tab
#
# 1) Refresh the per-row "time_taken" timestamps on the existing "data" sheet.
# 2) Add a new "metadata" worksheet (after "data") summarising the panel
#    query that produced this export.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:20:19.072581",
    "2021-10-05 14:20:19.072589",
    "2021-10-05 14:20:19.072593",
    "2021-10-05 14:20:19.072595",
    "2021-10-05 14:20:19.072598",
    "2021-10-05 14:20:19.072601",
    "2021-10-05 14:20:19.072604",
    "2021-10-05 14:20:19.072607",
    "2021-10-05 14:20:19.072609",
    "2021-10-05 14:20:19.072612",
    "2021-10-05 14:20:19.072615",
    "2021-10-05 14:20:19.072617",
    "2021-10-05 14:20:19.072620",
    "2021-10-05 14:20:19.072623",
    "2021-10-05 14:20:19.072625",
    "2021-10-05 14:20:19.072628",
    "2021-10-05 14:20:19.072631",
    "2021-10-05 14:20:19.072634",
    "2021-10-05 14:20:19.072637",
    "2021-10-05 14:20:19.072639",
    "2021-10-05 14:20:19.072642",
    "2021-10-05 14:20:19.072645",
    "2021-10-05 14:20:19.072647",
    "2021-10-05 14:20:19.072650",
    "2021-10-05 14:20:19.072653",
    "2021-10-05 14:20:19.072655",
    "2021-10-05 14:20:19.072658",
    "2021-10-05 14:20:19.072661",
    "2021-10-05 14:20:19.072663",
    "2021-10-05 14:20:19.072666",
    "2021-10-05 14:20:19.072669",
    "2021-10-05 14:20:19.072671",
    "2021-10-05 14:20:19.072674",
    "2021-10-05 14:20:19.072677",
    "2021-10-05 14:20:19.072680",
    "2021-10-05 14:20:19.072682",
    "2021-10-05 14:20:19.072685",
    "2021-10-05 14:20:19.072687",
    "2021-10-05 14:20:19.072690",
    "2021-10-05 14:20:19.072692",
    "2021-10-05 14:20:19.072695",
    "2021-10-05 14:20:19.072698",
    "2021-10-05 14:20:19.072701",
    "2021-10-05 14:20:19.072704"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the "metadata" worksheet, placed after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - reuse "data"'s header formatting (bold + border + centered)
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($c = 2; $c -le 7; $c++) {
    $meta.Cells.Item(1, $c).Value = $headers[$c - 2]
}
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$meta.Cells.Item(2, 2).Value = "Familial Neural Tube Defects"
$meta.Cells.Item(2, 3).Value = 11

# "1.10" must stay literal text (not collapse to the number 1.1) - force the
# cell to Text format for the assignment, then strip the formatting back off
# so the cell keeps the workbook's default (unstyled) appearance.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.10"
$meta.Cells.Item(2, 4).ClearFormats()

$meta.Cells.Item(2, 5).Value = "2021-01-15T12:07:16.183096Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:19.069145"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/11/?format=json"

# Index column (A2) - reuse "data"'s index-column formatting (bordered, centered)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Cells.Item(2, 1).Value = 0

$excel.CutCopyMode = $false

Write-Host "done"
